$wb = $excel.ActiveWorkbook

# --- Rename worksheets (tab renaming as part of the output-process refactor) ---
$renames = @{
    "PR_RSV_Uncut_Sheet"              = "PR_RSV_Uncut_Sheet_1"
    "PR_ADE_Uncut_Sheet"              = "PR_ADE_Uncut_Sheet_1"
    "PR_FLU_Buffer"                   = "PR_FLU_Buffer_1"
    "PR_FLU_Cassette"                 = "PR_FLU_Cassette_1"
    "PR_FLU_Extraction_Tubes_Dropper" = "PR_FLU_Ext_Tubes_Dropper_1"
    "PR_FLU_Strerile_Swabs"           = "PR_FLU_Strerile_Swabs_1"
    "PR_FSVA_Sterile_Swabs"           = "PR_FSVA_Sterile_Swabs_1"
    "PR_FSVA_Extraction_Tube"         = "PR_FSVA_Extraction_Tube_1"
    "PR_FSV_Buffer"                   = "PR_FSV_Buffer_1"
    "PR_FSV_Cassette"                 = "PR_FSV_Cassette_1"
    "PR_FSVA_Buffer"                  = "PR_FSVA_Buffer_1"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# --- Re-stamp each renamed sheet's print area so the workbook-level
#     Print_Area defined name text picks up the sheet's new name ---
$printAreas = @{
    "PR_RSV_Uncut_Sheet_1"       = '$A$1:$G$42'
    "PR_ADE_Uncut_Sheet_1"       = '$A$1:$G$42'
    "PR_FLU_Buffer_1"            = '$A$1:$G$42'
    "PR_FLU_Cassette_1"          = '$A$1:$H$42'
    "PR_FLU_Ext_Tubes_Dropper_1" = '$A$1:$G$42'
    "PR_FLU_Strerile_Swabs_1"    = '$A$1:$G$42'
    "PR_FSVA_Sterile_Swabs_1"    = '$A$1:$G$42'
    "PR_FSVA_Extraction_Tube_1"  = '$A$1:$G$42'
    "PR_FSV_Buffer_1"            = '$A$1:$G$42'
    "PR_FSV_Cassette_1"          = '$A$1:$H$42'
    "PR_FSVA_Buffer_1"           = '$A$1:$G$42'
}

foreach ($name in $printAreas.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.PageSetup.PrintArea = $printAreas[$name]
}

# --- Move the active tab: the newly-selected sheet becomes
#     PR_FLU_Ext_Tubes_Dropper_1 (activeTab index 6) ---
$wb.Worksheets.Item("PR_FLU_Ext_Tubes_Dropper_1").Activate()
